$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 25, shifting existing rows 25-45 down to 26-46
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with its data (copy format/values from row 26, which now holds what used to be row 25)
$ws.Range("A25:R25").Value = $ws.Range("A26:R26").Value2

$ws.Range("D25").Value = 44483
$ws.Range("J25").Value = 20
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 25000
$ws.Range("N25").Value = "$/malla 25 kilos"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 1000
